$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.487.02"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.910.34"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'244.46"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("D6").Value = "'0.9985"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.4827"
$ws.Range("E7").Value = "  +2.96%  "
$ws.Range("D8").Value = "'0.2892"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").Value = "'0.06698"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").Value = "'110.42"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").Value = "'19.18"
$ws.Range("E11").Value = "  +6.15%  "
$ws.Range("D12").Value = "1.913.23"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'0.07543"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "'5.254"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").Value = "'0.6670"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "'273.59"
$ws.Range("E16").Value = "  -4.40%  "
$ws.Range("D17").Value = "30.477.56"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "'0.9987"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "'0.000007527"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").Value = "'12.82"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").Value = "2.162.29"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "'5.478"
$ws.Range("E22").Value = "  +5.20%  "
$ws.Range("D23").Value = "'0.9993"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'6.414"
$ws.Range("E24").Value = "  +3.98%  "
$ws.Range("D25").Value = "'9.408"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("D26").Value = "'163.32"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").Value = "'20.15"
$ws.Range("E27").Value = "  -4.88%  "
$ws.Range("D28").Value = "'2.096"
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("D29").Value = "'0.1046"
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("D30").Value = "'1.398"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").Value = "'4.123"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'4.041"
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("D33").Value = "'0.04978"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").Value = "'0.7266"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").Value = "'1.130"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").Value = "'0.9988"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02031"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "'2.709"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").Value = "'2.667"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "'110.74"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("D42").Value = "'0.4411"
$ws.Range("E42").Value = "  +5.35%  "
$ws.Range("D43").Value = "'0.8666"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").Value = "'5.854"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "'0.9977"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "'67.60"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "'7.359"
$ws.Range("E47").Value = "  +3.57%  "
$ws.Range("D48").Value = "'9.288"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "'0.1241"
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("D50").Value = "'47.69"
$ws.Range("E50").Value = "  -9.21%  "
$ws.Range("D51").Value = "'1.462"
$ws.Range("E51").Value = "  +6.88%  "
